# Reorder the "Recorded By" values in column G of the "Session Analysis
# Results" sheet so that entries are listed with "System" (capitalized)
# first, matching the canonical ordering produced by the upstream sync.
#
# Only two distinct source strings are affected in this workbook:
#   "dnasr281@gmail.com, System"            -> "System, dnasr281@gmail.com"
#   "system, System, backup@backdoor.com"    -> "System, backup@backdoor.com, system"
# All other values in column G are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")
if ($null -eq $ws) {
    $ws = $wb.ActiveSheet
}

$replacements = @{
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "system, System, backup@backdoor.com" = "System, backup@backdoor.com, system"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value2

    if ($null -ne $val -and $replacements.ContainsKey($val)) {
        $cell.Value2 = $replacements[$val]
    }
}
